# cryptos.xlsx refresh: updates the "Price" (D) and "Volume(1h)" (E) columns
# on Sheet1 to the latest scraped coinranking.com figures, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.965.46"
$ws.Range("E2").Value = "  +2.52%  "

$ws.Range("D3").Value = "3.729.83"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.75"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.95"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "3.729.36"
$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("E11").Value = "  +3.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.75"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").Value = "4.357.58"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").Value = "3.721.07"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "69.039.75"
$ws.Range("E17").Value = "  +2.54%  "

$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.93"
$ws.Range("E20").Value = "  -2.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.70"
$ws.Range("E21").Value = "  +15.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.78"
$ws.Range("E22").Value = "  +0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.721"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("E24").Value = "  +4.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.62"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.49"
$ws.Range("E31").Value = "  +5.94%  "

$ws.Range("E32").Value = "  +4.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.38"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").Value = "3.880.25"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "3.668.92"
$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.96"
$ws.Range("E42").Value = "  +4.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "430.40"
$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.59"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("E45").Value = "  +1.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.43"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.06"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.43"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").Value = "2.761.36"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("E51").Value = "  +0.54%  "

# Drop the temporary "@" text-number-format override again so the cells end up
# with no explicit style, exactly like the rest of the data rows.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
